$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '73.065.37'
$ws.Range('E2').Value = '  +6.01%  '
$ws.Range('D3').Value = '2.589.89'
$ws.Range('E3').Value = '  +6.12%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''580.87'
$ws.Range('E5').Value = '  +3.69%  '
$ws.Range('D6').Value = '''182.56'
$ws.Range('E6').Value = '  +12.54%  '
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '''0.532'
$ws.Range('E8').Value = '  +2.79%  '
$ws.Range('D9').Value = '''0.193'
$ws.Range('E9').Value = '  +14.40%  '
$ws.Range('D10').Value = '2.588.46'
$ws.Range('E10').Value = '  +6.01%  '
$ws.Range('D11').Value = '''0.162'
$ws.Range('E11').Value = '  -0.14%  '
$ws.Range('D12').Value = '''0.358'
$ws.Range('E12').Value = '  +7.91%  '
$ws.Range('D13').Value = '''4.71'
$ws.Range('E13').Value = '  +2.32%  '
$ws.Range('D14').Value = '3.112.53'
$ws.Range('E14').Value = '  +7.70%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '73.320.17'
$ws.Range('E15').Value = '  +6.54%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '''0.0000184'
$ws.Range('E16').Value = '  +4.85%  '
$ws.Range('D17').Value = '''25.76'
$ws.Range('E17').Value = '  +10.99%  '
$ws.Range('D18').Value = '2.620.68'
$ws.Range('E18').Value = '  +7.44%  '
$ws.Range('D19').Value = '''8.78'
$ws.Range('E19').Value = '  +26.44%  '
$ws.Range('D20').Value = '''11.67'
$ws.Range('E20').Value = '  +11.01%  '
$ws.Range('D21').Value = '''368.59'
$ws.Range('E21').Value = '  +8.69%  '
$ws.Range('D22').Value = '''2.21'
$ws.Range('E22').Value = '  +14.96%  '
$ws.Range('D23').Value = '''4.04'
$ws.Range('E23').Value = '  +5.38%  '
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('D25').Value = '''69.26'
$ws.Range('E25').Value = '  +3.22%  '
$ws.Range('D26').Value = '''4.09'
$ws.Range('E26').Value = '  +10.39%  '
$ws.Range('D27').Value = '''9.16'
$ws.Range('E27').Value = '  +11.50%  '
$ws.Range('E28').Value = '  +6.70%  '
$ws.Range('E29').Value = '  -0.87%  '
$ws.Range('D30').Value = '0.0₃0923'
$ws.Range('E30').Value = '  +12.94%  '
$ws.Range('D31').Value = '''7.84'
$ws.Range('E31').Value = '  +9.91%  '
$ws.Range('D32').Value = '''1.36'
$ws.Range('E32').Value = '  +17.16%  '
$ws.Range('D33').Value = '''499.50'
$ws.Range('E33').Value = '  +16.59%  '
$ws.Range('D34').Value = '''1.72'
$ws.Range('E34').Value = '  +6.82%  '
$ws.Range('E35').Value = '  -0.24%  '
$ws.Range('D36').Value = '''0.120'
$ws.Range('E36').Value = '  +12.74%  '
$ws.Range('D37').Value = '''159.82'
$ws.Range('E37').Value = '  +0.24%  '
$ws.Range('D38').Value = '''18.96'
$ws.Range('E38').Value = '  +5.37%  '
$ws.Range('D39').Value = '''19.24'
$ws.Range('E39').Value = '  +1.06%  '
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('D41').Value = '''4.81'
$ws.Range('E41').Value = '  +10.82%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '''1.64'
$ws.Range('E42').Value = '  +9.21%  '
$ws.Range('B43').Value = 'PolygonEcosystemToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D43').Value = '''0.321'
$ws.Range('E43').Value = '  +7.60%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').Value = '''0.0895'
$ws.Range('E44').Value = '  +24.82%  '
$ws.Range('D45').Value = '''155.74'
$ws.Range('E45').Value = '  +19.41%  '
$ws.Range('D46').Value = '''2.33'
$ws.Range('E46').Value = '  +14.53%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').Value = '''38.60'
$ws.Range('E47').Value = '  +3.16%  '
$ws.Range('B48').Value = 'ImmutableX'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D48').Value = '''1.16'
$ws.Range('E48').Value = '  +7.43%  '
$ws.Range('D49').Value = '''3.58'
$ws.Range('E49').Value = '  +6.88%  '
$ws.Range('D50').Value = '''0.519'
$ws.Range('E50').Value = '  +7.69%  '
$ws.Range('D51').Value = '''19.95'
$ws.Range('E51').Value = '  +18.18%  '
